# edit.ps1 - reproduces:
#   1. Table (slide 16, shape 3) tableStyleId change
#        {5045F0DA-FFD9-4E95-BA13-E6A5609EB5B4} -> {5FCC5CF3-892C-47E2-87CB-E4B8DC153804}
#   2. The presentation's theme colour scheme is swapped from the custom
#      "Integral" palette back to the stock "Office Theme" palette.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style change
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
$tableShape = $tableSlide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{5FCC5CF3-892C-47E2-87CB-E4B8DC153804}")

# ---------------------------------------------------------------------------
# 2) Theme colour scheme swap (Integral -> Office Theme)
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$colors = $master.Theme.ThemeColorScheme

# Index : Office Theme RGB value (COM packed 0x00BBGGRR)
#   1 dk1      000000
#   2 lt1      FFFFFF
#   3 dk2      44546A
#   4 lt2      E7E6E6
#   5 accent1  5B9BD5
#   6 accent2  ED7D31
#   7 accent3  A5A5A5
#   8 accent4  FFC000
#   9 accent5  4472C4
#  10 accent6  70AD47
#  11 hlink    0563C1
#  12 folHlink 954F72
$colors.Item(1).RGB  = 0
$colors.Item(2).RGB  = 16777215
$colors.Item(3).RGB  = 6968388
$colors.Item(4).RGB  = 15132391
$colors.Item(5).RGB  = 13998939
$colors.Item(6).RGB  = 3243501
$colors.Item(7).RGB  = 10855845
$colors.Item(8).RGB  = 49407
$colors.Item(9).RGB  = 12874308
$colors.Item(10).RGB = 4697456
$colors.Item(11).RGB = 12673797
$colors.Item(12).RGB = 7491477
